$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on D-column cells whose new values would
# otherwise be auto-parsed as numbers (losing trailing zeros / exact text),
# matching the original inline-string "Price" column formatting.
$u = $ws.Range("D5")
$u = $excel.Union($u, $ws.Range("D6"))
$u = $excel.Union($u, $ws.Range("D7"))
$u = $excel.Union($u, $ws.Range("D9"))
$u = $excel.Union($u, $ws.Range("D10"))
$u = $excel.Union($u, $ws.Range("D11"))
$u = $excel.Union($u, $ws.Range("D14"))
$u = $excel.Union($u, $ws.Range("D15"))
$u = $excel.Union($u, $ws.Range("D16"))
$u = $excel.Union($u, $ws.Range("D20"))
$u = $excel.Union($u, $ws.Range("D21"))
$u = $excel.Union($u, $ws.Range("D22"))
$u = $excel.Union($u, $ws.Range("D25"))
$u = $excel.Union($u, $ws.Range("D26"))
$u = $excel.Union($u, $ws.Range("D27"))
$u = $excel.Union($u, $ws.Range("D28"))
$u = $excel.Union($u, $ws.Range("D29"))
$u = $excel.Union($u, $ws.Range("D33"))
$u = $excel.Union($u, $ws.Range("D34"))
$u = $excel.Union($u, $ws.Range("D35"))
$u = $excel.Union($u, $ws.Range("D36"))
$u = $excel.Union($u, $ws.Range("D38"))
$u = $excel.Union($u, $ws.Range("D40"))
$u = $excel.Union($u, $ws.Range("D42"))
$u = $excel.Union($u, $ws.Range("D43"))
$u = $excel.Union($u, $ws.Range("D44"))
$u = $excel.Union($u, $ws.Range("D45"))
$u = $excel.Union($u, $ws.Range("D46"))
$u = $excel.Union($u, $ws.Range("D48"))
$u.NumberFormat = "@"

$ws.Range("D2").Value = "43.922.40"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.355.27"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.672"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").Value = "237.35"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("D7").Value = "72.68"
$ws.Range("E7").Value = "  +9.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +19.24%  "
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "28.15"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "2.702.81"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "16.67"
$ws.Range("E14").Value = "  +7.51%  "
$ws.Range("D15").Value = "6.66"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").Value = "0.893"
$ws.Range("E16").Value = "  +4.73%  "
$ws.Range("D17").Value = "2.359.25"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "43.868.25"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "77.77"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").Value = "6.44"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "254.53"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "10.55"
$ws.Range("E26").Value = "  +5.78%  "
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "22.40"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "172.34"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  +5.96%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("D34").Value = "0.0712"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").Value = "4.00"
$ws.Range("E36").Value = "  +10.38%  "
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").Value = "6.42"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +6.16%  "
$ws.Range("D40").Value = "19.71"
$ws.Range("E40").Value = "  +7.99%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "8.83"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").Value = "1.23"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "0.0979"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").Value = "1.16"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").Value = "4.45"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("E47").Value = "  +10.91%  "
$ws.Range("D48").Value = "97.57"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "1.434.63"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +1.38%  "
